$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells need a leading apostrophe to force Excel to
# keep the numeric-looking text as a string instead of auto-converting
# it to a floating point number (which would lose the exact formatting).
$ws.Range('D2').Value = '''26.957.25'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '''1.553.20'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = '''207.15'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('E6').Value = '  +0.98%  '
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').Value = '''21.57'
$ws.Range('E9').Value = '  +0.80%  '
$ws.Range('D10').Value = '''0.0588'
$ws.Range('E10').Value = '  +1.28%  '
$ws.Range('D11').Value = '''0.0861'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '''1.776.02'
$ws.Range('E12').Value = '  +0.60%  '
$ws.Range('D13').Value = '''1.555.08'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '''3.71'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('D16').Value = '''61.91'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = '''26.946.55'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '''215.30'
$ws.Range('E18').Value = '  +1.26%  '
$ws.Range('D19').Value = '''0.0₃0689'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').Value = '''7.25'
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').Value = '''9.18'
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').Value = '''152.42'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('D26').Value = '''6.65'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').Value = '''14.88'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('E28').Value = '  +0.58%  '
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('E30').Value = '  +0.71%  '
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').Value = '''3.21'
$ws.Range('E32').Value = '  +0.64%  '
$ws.Range('D33').Value = '''1.403.61'
$ws.Range('E33').Value = '  +5.37%  '
$ws.Range('E34').Value = '  +3.21%  '
$ws.Range('D35').Value = '''1.55'
$ws.Range('E35').Value = '  +3.04%  '
$ws.Range('D36').Value = '''0.952'
$ws.Range('E36').Value = '  +2.24%  '
$ws.Range('E37').Value = '  +0.63%  '
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').Value = '''0.522'
$ws.Range('E39').Value = '  -0.06%  '
$ws.Range('D40').Value = '''0.808'
$ws.Range('E40').Value = '  +1.03%  '
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('D42').Value = '''0.989'
$ws.Range('E42').Value = '  -0.51%  '
$ws.Range('E43').Value = '  +3.47%  '
$ws.Range('E44').Value = '  -4.07%  '
$ws.Range('D45').Value = '''63.75'
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('D46').Value = '''1.73'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').Value = '''1.690.11'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = '''86.24'
$ws.Range('E48').Value = '  +0.43%  '
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('D50').Value = '''0.0956'
$ws.Range('E50').Value = '  +0.47%  '
$ws.Range('E51').Value = '  +0.56%  '
